$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert the "In my own words" answer for Problem #1 right after the
#    cat/parrot/seed problem statement paragraph, before the blank
#    separator paragraph that precedes "Socks in the dark:".
# ---------------------------------------------------------------------

# Locate the paragraph that states the cat/parrot/seed riddle.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "*without leaving the wrong ones alone together.*") {
        $target = $para
        break
    }
}

# The paragraph immediately following it is the blank "<w:p/>" separator
# that comes right before "Socks in the dark:". We insert our two new
# paragraphs right before that separator (i.e. right after the riddle).
$sepRange = $target.Next().Range
$sepRange.Collapse(1)   # wdCollapseStart
$sepRange.InsertBefore(" `rA. The problem is the man needs to figure out the best way to get everyone safely across the river without risking the animals eating each other or the food and he can only take one item at a time.`r")

# ---------------------------------------------------------------------
# 2. Move the "_GoBack" bookmark from the end of the document onto the
#    newly-typed answer paragraph (this mirrors Word's own behaviour of
#    stamping _GoBack at the most-recently-edited location).
# ---------------------------------------------------------------------

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Re-find the answer paragraph (now that text has shifted) by content.
$answerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "A. The problem is the man needs to figure out*") {
        $answerPara = $para
        break
    }
}

$bmRange = $d.Range($answerPara.Range.Start, $answerPara.Range.End - 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Output "Edit applied."
